$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("8x8")

# Update the reference value in A12; dependent formulas in C11:E13 recalc.
$ws.Range("A12").Value = 24

# Update the active selection to match the saved view state.
$ws.Range("A12").Select()
